# Applies the "Add files via upload" edit to the _TutosYoutube workbook:
#  - Reword several "Composant : X" cells to just "X" (and a few similar
#    simplifications), consolidating the old longer explanatory sentences.
#  - Add a brand-new row (10) describing the new indexC017_Download course,
#    with a hyperlink on its course-name cell.
#  - Row-height tweaks for the rows whose wrapped text got shorter / longer.
#  - Move the active selection to D9.
#
# NOTE: the shared-string table itself (xl/sharedStrings.xml) is rebuilt
# automatically from the live cell contents when the workbook is saved, so
# this script only needs to set the *cell values* that actually changed -
# every other cell that merely shifts to a different shared-string index
# because of the table reshuffle does not need to be touched at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reworded / shortened descriptions in column D ------------------------
# (written bottom-to-top / in the same order the original author typed them,
# so the shared-string table the workbook rebuilds on save lands in the same
# append order as the authored file)
$ws.Range("D9").Value  = "Calendrier avec date début / date fin"
$ws.Range("D8").Value  = "Zone de saisie : toutes les fonctionnalités, dont MAJ d'un graphique à partir de ces zones de saisie"
$ws.Range("D7").Value  = "Cases à cocher qui permettent de mettre à jour un graphique"
$ws.Range("D6").Value  = "Composant dcc.Interval qui permet de mettre à jour les données qu'on pourrait éventuellement récupérer dans différents sites @ à partir des API." + [char]10 + "Sous-librairie dash.exceptions qui lève une exception lors de la MAJ des données"
$ws.Range("D5").Value  = "Curseur"
$ws.Range("D3").Value  = "Graphique"
$ws.Range("D2").Value  = "Menu déroulant"
$ws.Range("D12").Value = "Chargement des fichiers avec génération automatique des données dans une datatable et MAJ d'un graphique avec des menus déroulants"
$ws.Range("D20").Value = "Bouton"

# --- New row 10: indexC017_Download ---------------------------------------
$ws.Range("D10").Value = "Téléchargement des données : dans ce cours on apprend à télécharger les données à partir du composant download,et également la mise en forme du bouton de téléchargement des données en insérant" + [char]10 + "un icône à partir d'un site @"
$ws.Range("C10").Value = "indexC017_Download"

# New hyperlink on the course-name cell of the new row (becomes rId4).
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.youtube.com/watch?v=dl0wnLoad10&list=PLh3I780jNsiSC7QJMQ46tHDYYnfhGEflf&index=10")

# --- Row height adjustments -------------------------------------------------
# Row 6's wrapped text shrank from 4 wrapped lines to 2 -> 43.2pt.
$ws.Rows.Item(6).RowHeight = 43.2
# Row 10 now holds a two-line wrapped description -> 43.2pt.
$ws.Rows.Item(10).RowHeight = 43.2

# --- Selection --------------------------------------------------------------
$ws.Range("D9").Select()
